$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.813.76"
$ws.Range("E2").Value = "  +0.66%  "

# Row 3
$ws.Range("D3").Value = "2.290.01"

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'96.29"
$ws.Range("E5").Value = "  +3.79%  "

# Row 6
$ws.Range("D6").Value = "'269.87"
$ws.Range("E6").Value = "  +1.16%  "

# Row 7
$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("E9").Value = "  -0.98%  "

# Row 10
$ws.Range("D10").Value = "'45.59"
$ws.Range("E10").Value = "  +2.96%  "

# Row 11
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").Value = "'7.90"
$ws.Range("E12").Value = "  -1.56%  "

# Row 13
$ws.Range("E13").Value = "  +1.88%  "

# Row 14
$ws.Range("D14").Value = "'15.77"
$ws.Range("E14").Value = "  +3.15%  "

# Row 15
$ws.Range("D15").Value = "2.633.61"
$ws.Range("E15").Value = "  -0.79%  "

# Row 16
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
$ws.Range("D17").Value = "2.292.49"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18
$ws.Range("D18").Value = "43.684.05"
$ws.Range("E18").Value = "  +0.48%  "

# Row 19
$ws.Range("E19").Value = "  +4.26%  "

# Row 20
$ws.Range("E20").Value = "  -2.21%  "

# Row 21
$ws.Range("D21").Value = "'72.14"
$ws.Range("E21").Value = "  +1.30%  "

# Row 22
$ws.Range("D22").Value = "'2.50"
$ws.Range("E22").Value = "  +11.32%  "

# Row 23
$ws.Range("E23").Value = "  -1.47%  "

# Row 24
$ws.Range("E24").Value = "  -4.29%  "

# Row 25
$ws.Range("D25").Value = "'2.76"
$ws.Range("E25").Value = "  +11.21%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("E27").Value = "  +0.71%  "

# Row 28
$ws.Range("E28").Value = "  -1.71%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.28"
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'38.65"
$ws.Range("E30").Value = "  +0.43%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'174.70"
$ws.Range("E31").Value = "  +2.04%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'21.84"
$ws.Range("E32").Value = "  -3.06%  "

# Row 33
$ws.Range("D33").Value = "'0.0898"

# Row 34
$ws.Range("E34").Value = "  -0.65%  "

# Row 35
$ws.Range("E35").Value = "  +0.67%  "

# Row 36
$ws.Range("E36").Value = "  +4.14%  "

# Row 37
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("E38").Value = "  -1.69%  "

# Row 39
$ws.Range("D39").Value = "'3.51"
$ws.Range("E39").Value = "  +3.47%  "

# Row 40
$ws.Range("E40").Value = "  +1.30%  "

# Row 41
$ws.Range("E41").Value = "  -0.38%  "

# Row 42
$ws.Range("D42").Value = "'12.30"
$ws.Range("E42").Value = "  +2.69%  "

# Row 43
$ws.Range("D43").Value = "'1.33"
$ws.Range("E43").Value = "  -0.96%  "

# Row 44
$ws.Range("D44").Value = "'64.13"
$ws.Range("E44").Value = "  +4.69%  "

# Row 45
$ws.Range("D45").Value = "'8.76"
$ws.Range("E45").Value = "  -2.41%  "

# Row 46
$ws.Range("D46").Value = "'5.20"
$ws.Range("E46").Value = "  -3.03%  "

# Row 47
$ws.Range("E47").Value = "  +0.06%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'97.56"
$ws.Range("E48").Value = "  -2.11%  "

# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.20"
$ws.Range("E49").Value = "  -0.13%  "

# Row 50
$ws.Range("D50").Value = "'1.53"
$ws.Range("E50").Value = "  +13.05%  "

# Row 51
$ws.Range("D51").Value = "'0.432"
$ws.Range("E51").Value = "  +4.30%  "
